$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.249.35"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "3.821.17"
$ws.Range("E3").Value = "  -0.65%  "
$c = $ws.Range("D4")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "
$c = $ws.Range("D5")
$c.Value = "'705.46"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").Value = "3.819.28"
$ws.Range("E7").Value = "  -0.66%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  -0.16%  "
$c = $ws.Range("D11")
$c.Value = "'7.58"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +4.04%  "
$c = $ws.Range("D12")
$c.Value = "'0.463"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("E13").Value = "  -1.08%  "
$c = $ws.Range("D14")
$c.Value = "'36.06"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").Value = "4.460.77"
$ws.Range("E15").Value = "  -0.72%  "
$ws.Range("D16").Value = "3.795.86"
$ws.Range("E16").Value = "  -2.18%  "
$ws.Range("D17").Value = "71.224.32"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c = $ws.Range("D18")
$c.Value = "'17.56"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D19")
$c.Value = "'7.17"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("E20").Value = "  -0.59%  "
$c = $ws.Range("D21")
$c.Value = "'511.82"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +3.68%  "
$c = $ws.Range("D22")
$c.Value = "'10.72"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("E23").Value = "  +0.90%  "
$c = $ws.Range("D24")
$c.Value = "'84.28"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.32%  "
$c = $ws.Range("D25")
$c.Value = "'0.0000145"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.81%  "
$ws.Range("D26").Value = "3.968.10"
$ws.Range("E26").Value = "  -0.76%  "
$c = $ws.Range("D27")
$c.Value = "'12.07"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.64%  "
$c = $ws.Range("D28")
$c.Value = "'10.42"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.34%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("E30").Value = "  -3.69%  "
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("E33").Value = "  -0.60%  "
$c = $ws.Range("D34")
$c.Value = "'29.18"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.91%  "
$c = $ws.Range("D35")
$c.Value = "'0.174"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -4.58%  "
$c = $ws.Range("D36")
$c.Value = "'9.18"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("D37").Value = "3.778.91"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("E41").Value = "  +0.46%  "
$c = $ws.Range("D42")
$c.Value = "'1.03"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -1.45%  "
$c = $ws.Range("D43")
$c.Value = "'3.32"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.52%  "
$c = $ws.Range("D45")
$c.Value = "'172.01"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +5.16%  "
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("E47").Value = "  -0.53%  "
$c = $ws.Range("D48")
$c.Value = "'49.69"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +2.03%  "
$c = $ws.Range("D49")
$c.Value = "'430.21"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +6.15%  "
$c = $ws.Range("D50")
$c.Value = "'8.68"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.62%  "
$c = $ws.Range("D51")
$c.Value = "'0.294"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.46%  "
